# Update the Cxcl13-Cxcr5 NATMI sending/target-cluster pair table with the new TPM-based run.
# The "Resolving-Mac" sending cluster is renamed to "MuSCs", the target-cluster set grows from
# {FAPs, MuSCs, Neutrophils} to {ECs, FAPs, Inflammatory-Mac, MuSCs, Neutrophils}, and every
# metric column (E:T) is refreshed for the full 2 (sending) x 5 (target) cluster grid -- rows 2-11.
# The header row (row 1) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: A..T = Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, then the 16 numeric NATMI metric columns.
$rows = @(
    @("FAPs", "Cxcl13", "Cxcr5", "ECs", 3, 1, 7.767740666666666, 23.303222, 0.957755623847744, 0.9714348434930534, 2, 1, 0.8542635000000001, 1.708527, 0.2602038345382578, 0.1959844976819479, 6.635697328999, 39.814183973994, 0.2492116858757643, 0.1903861698327277),
    @("FAPs", "Cxcl13", "Cxcr5", "FAPs", 3, 1, 7.767740666666666, 23.303222, 0.957755623847744, 0.9714348434930534, 1, 0.3333333333333333, 0.1625146666666667, 0.487544, 0.0495010490737032, 0.05592599118295911, 1.262371785196444, 11.361346066768, 0.0474099081367024, 0.05432845649201177),
    @("FAPs", "Cxcl13", "Cxcr5", "Inflammatory-Mac", 3, 1, 7.767740666666666, 23.303222, 0.957755623847744, 0.9714348434930534, 1, 0.3333333333333333, 0.130141, 0.390423, 0.03964021315512534, 0.04478527734035173, 1.010901538100667, 9.098113842905999, 0.03796563707984462, 0.04350597888391758),
    @("FAPs", "Cxcl13", "Cxcr5", "MuSCs", 3, 1, 7.767740666666666, 23.303222, 0.957755623847744, 0.9714348434930534, 1, 0.5, 0.2772375, 0.554475, 0.08444497579236412, 0.06360362133709214, 2.153509003075, 12.92105401845, 0.08087765047082335, 0.06178677393918954),
    @("FAPs", "Cxcl13", "Cxcr5", "Neutrophils", 3, 1, 7.767740666666666, 23.303222, 0.957755623847744, 0.9714348434930534, 3, 1, 1.858898333333333, 5.576695, 0.5662099274405495, 0.6397006124576493, 14.43944017903222, 129.95496161129, 0.5422907422846094, 0.621427464345207),
    @("MuSCs", "Cxcl13", "Cxcr5", "ECs", 1, 0.5, 0.342617, 0.685234, 0.04224437615225601, 0.02856515650694651, 2, 1, 0.8542635000000001, 1.708527, 0.2602038345382578, 0.1959844976819479, 0.2926851975795, 1.170740790318, 0.01099214866249355, 0.005598327849220136),
    @("MuSCs", "Cxcl13", "Cxcr5", "FAPs", 1, 0.5, 0.342617, 0.685234, 0.04224437615225601, 0.02856515650694651, 1, 0.3333333333333333, 0.1625146666666667, 0.487544, 0.0495010490737032, 0.05592599118295911, 0.05568028754933334, 0.334081725296, 0.002091140937000802, 0.001597534690947338),
    @("MuSCs", "Cxcl13", "Cxcr5", "Inflammatory-Mac", 1, 0.5, 0.342617, 0.685234, 0.04224437615225601, 0.02856515650694651, 1, 0.3333333333333333, 0.130141, 0.390423, 0.03964021315512534, 0.04478527734035173, 0.044588518997, 0.267531113982, 0.001674576075280722, 0.001279298456434152),
    @("MuSCs", "Cxcl13", "Cxcr5", "MuSCs", 1, 0.5, 0.342617, 0.685234, 0.04224437615225601, 0.02856515650694651, 1, 0.5, 0.2772375, 0.554475, 0.08444497579236412, 0.06360362133709214, 0.09498628053750001, 0.37994512215, 0.003567325321540783, 0.0018168473979026),
    @("MuSCs", "Cxcl13", "Cxcr5", "Neutrophils", 1, 0.5, 0.342617, 0.685234, 0.04224437615225601, 0.02856515650694651, 3, 1, 1.858898333333333, 5.576695, 0.5662099274405495, 0.6397006124576493, 0.6368901702716667, 3.82134102163, 0.02391918515594016, 0.01827314811244229)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $targetRow = $r + 2
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($targetRow, $c).Value = $rowData[$c - 1]
    }
}